$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value (cryptos.xlsx price/volume refresh)
$updates = [ordered]@{
    "D2" = "68.560.48"
    "E2" = "  +0.37%  "
    "D3" = "2.451.00"
    "E3" = "  -0.03%  "
    "E4" = "  -0.03%  "
    "D5" = "557.21"
    "E5" = "  -0.50%  "
    "D6" = "161.02"
    "E6" = "  -1.33%  "
    "E7" = "  -0.07%  "
    "E8" = "  +1.05%  "
    "D9" = "0.151"
    "E9" = "  +0.05%  "
    "E10" = "  +0.64%  "
    "D11" = "0.330"
    "E11" = "  -2.45%  "
    "D12" = "4.83"
    "E12" = "  +0.18%  "
    "D13" = "68.481.22"
    "E13" = "  +0.35%  "
    "D14" = "0.0000167"
    "E14" = "  -1.85%  "
    "D15" = "23.29"
    "E15" = "  -0.04%  "
    "D16" = "10.54"
    "E16" = "  -3.88%  "
    "D17" = "334.02"
    "E17" = "  -2.38%  "
    "D18" = "6.88"
    "E18" = "  -4.08%  "
    "D19" = "3.75"
    "E19" = "  -1.07%  "
    "B20" = "Dai"
    "C20" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D20" = "1.00"
    "E20" = "  +0.12%  "
    "B21" = "SuiNetwork"
    "C21" = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
    "D21" = "1.88"
    "E21" = "  +0.87%  "
    "D22" = "66.18"
    "E22" = "  -2.34%  "
    "D23" = "3.61"
    "E23" = "  -2.81%  "
    "D24" = "8.13"
    "E24" = "  -0.26%  "
    "D25" = "0.0₃0811"
    "E25" = "  -2.85%  "
    "D26" = "7.14"
    "E26" = "  -1.32%  "
    "D27" = "0.999"
    "E27" = "  -0.13%  "
    "D28" = "425.33"
    "E28" = "  -1.65%  "
    "D29" = "1.13"
    "E29" = "  -3.53%  "
    "D30" = "1.60"
    "E30" = "  -4.18%  "
    "D31" = "159.64"
    "E31" = "  +1.88%  "
    "D32" = "18.99"
    "E32" = "  -0.03%  "
    "E33" = "  -0.01%  "
    "E34" = "  -1.95%  "
    "D35" = "17.73"
    "E35" = "  -0.79%  "
    "D36" = "0.298"
    "E36" = "  -2.63%  "
    "D37" = "4.37"
    "E37" = "  -1.74%  "
    "D38" = "1.45"
    "E38" = "  -4.19%  "
    "D39" = "1.07"
    "E39" = "  -2.23%  "
    "D40" = "2.03"
    "E40" = "  -1.47%  "
    "D41" = "3.33"
    "E41" = "  -0.84%  "
    "D42" = "128.67"
    "E42" = "  -3.87%  "
    "D43" = "0.0716"
    "E43" = "  -0.20%  "
    "D44" = "0.478"
    "E44" = "  -1.18%  "
    "D45" = "0.560"
    "E45" = "  -0.07%  "
    "D46" = "0.0909"
    "E46" = "  +0.00%  "
    "E47" = "  +0.28%  "
    "D48" = "1.37"
    "E48" = "  -3.78%  "
    "D49" = "4.90"
    "E49" = "  -8.67%  "
    "D50" = "16.67"
    "E50" = "  -5.02%  "
    "D51" = "0.581"
    "E51" = "  -3.26%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text format so numeric-looking strings (e.g. "1.00", "0.0000167")
    # are not silently coerced into numbers / scientific notation.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

Write-Host "Updated $($updates.Count) cells"
